# Update row 2 result values on each of the six year sheets (2025, 2030,
# 2035, 2040, 2045, 2050) with the latest figures received from the server.

$wb = $excel.ActiveWorkbook

# Map: sheet index (1-based, in workbook order) -> column letter -> new value.
# Only the columns listed here change; all other columns on row 2 keep
# their existing values.
$updates = @{
    1 = @{ 'A' = 0;                 'E' = 29970.44631501978; 'G' = 8095.925712661508;
           'I' = 13738.00533864;    'L' = 54416.76749186649; 'M' = 10518.579755365;
           'N' = 7532.916175655135; 'O' = 6996.064443878286 }
    2 = @{ 'A' = 0;                 'B' = 3098.910837094656; 'E' = 47468.9700841482;
           'G' = 8095.925712661508; 'I' = 22151.04139999618; 'L' = 79845.75677560513;
           'M' = 16283.86611051475; 'N' = 9392.991424751412; 'O' = 8237.889785501595 }
    3 = @{ 'A' = 1495.31406888776;  'B' = 5882.712195850656; 'E' = 60951.26884845589;
           'G' = 8095.925712661508; 'I' = 39923.28950426233; 'L' = 79845.75677560513;
           'M' = 20960.87118907375; 'N' = 13774.23956835204; 'O' = 12885.20435109773 }
    4 = @{ 'A' = 1495.31406888776;  'B' = 5882.712195850656; 'E' = 60951.26884845589;
           'G' = 8095.925712661508; 'I' = 39923.28950426233; 'L' = 79845.75677560513;
           'M' = 20960.87118907375; 'N' = 13774.23956835204; 'O' = 12885.20435109773 }
    5 = @{ 'A' = 1495.31406888776;  'B' = 5882.712195850656; 'E' = 60951.26884845589;
           'G' = 8095.925712661508; 'I' = 39923.28950426233; 'L' = 79845.75677560513;
           'M' = 20960.87118907375; 'N' = 13774.23956835204; 'O' = 12885.20435109773 }
    6 = @{ 'A' = 1495.31406888776;  'B' = 5882.712195850656; 'E' = 60951.26884845589;
           'G' = 8095.925712661508; 'I' = 39923.28950426233; 'L' = 79845.75677560513;
           'M' = 20960.87118907375; 'N' = 13774.23956835204; 'O' = 12885.20435109773 }
}

foreach ($sheetIndex in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $cols = $updates[$sheetIndex]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col`2").Value = $cols[$col]
    }
}
